# Daily auto push: insert the 2026/01/06 09:00 data point into the
# time-series table on Sheet1. This shifts every following row down by
# one and extends the used range from A1:D616 to A1:D617.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 575 (pushes old rows 575..616 down to 576..617,
# exactly matching how Excel's Rows.Insert behaves interactively).
$ws.Rows.Item(575).Insert()

# Column A holds the date as plain text (e.g. "2026/01/06"), not a real
# date value, so force a Text number format before writing it -- otherwise
# Excel would silently coerce the "yyyy/mm/dd"-looking string into a date
# serial number.
$ws.Range("A575").NumberFormat = "@"
$ws.Range("A575").Value = "2026/01/06"
$ws.Range("B575").Value = "火"
$ws.Range("C575").Value = 9
$ws.Range("D575").Value = 201
